$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NH")

# The "county" (A) and "town" (B) columns had been entered swapped; fix by
# exchanging the values of columns A and B for every data row.
for ($r = 2; $r -le 235; $r++) {
    $a = $ws.Cells.Item($r, 1).Value2
    $b = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value = $b
    $ws.Cells.Item($r, 2).Value = $a
}

# Column widths now correspond to the new (swapped) contents.
$ws.Columns.Item(1).ColumnWidth = 12.7
$ws.Columns.Item(2).ColumnWidth = 14.35

# Make NH the active sheet/tab, with B6 selected.
$ws.Activate()
$ws.Range("B6").Select()
